$wb = $excel.ActiveWorkbook

# -----------------------------------------------------------------
# Sheet "login": rows for chrome/MozillaFirefox swap places, and a
# new row for the "ie" browser is appended (mirrors the chrome row).
# -----------------------------------------------------------------
$wsLogin = $wb.Worksheets.Item("login")

# Row 2 becomes the "chrome" entry (was row 3)
$wsLogin.Range("A2").Value = "chrome"
$wsLogin.Range("C2").Value = "vishnu"
$wsLogin.Range("D2").Value = "vishnu"

# Row 3 becomes the "MozillaFirefox" entry (was row 2)
$wsLogin.Range("A3").Value = "MozillaFirefox"
$wsLogin.Range("C3").Value = "admin"
$wsLogin.Range("D3").Value = "manager"

# Row 4 (new): "ie" browser, same url/credentials pattern as chrome row
$wsLogin.Range("A4").Value = "ie"
$wsLogin.Hyperlinks.Add($wsLogin.Range("B4"), "http://localhost:9090/login.do")
$wsLogin.Range("B4").Style = $wsLogin.Range("B2").Style
$wsLogin.Range("C4").Value = "vishnu"
$wsLogin.Range("D4").Value = "vishnu"

$wsLogin.Rows.Item(2).Select()

# -----------------------------------------------------------------
# Sheet "CreateUser": add a new "status" column H (pass/fail per row)
# -----------------------------------------------------------------
$wsUser = $wb.Worksheets.Item("CreateUser")

$wsUser.Range("H2").Value = "pass"
$wsUser.Range("H3").Value = "fail"
$wsUser.Range("H4").Value = "pass"
$wsUser.Range("H5").Value = "fail"
$wsUser.Range("H6").Value = "fail"
$wsUser.Range("H7").Value = "pass"

$wsUser.Range("E13").Select()

# -----------------------------------------------------------------
# Sheet "CreateCustomer": selection only changes
# -----------------------------------------------------------------
$wsCustomer = $wb.Worksheets.Item("CreateCustomer")
$wsCustomer.Range("A8").Select()

# -----------------------------------------------------------------
# Sheet "CreateProject": selection only changes
# -----------------------------------------------------------------
$wsProject = $wb.Worksheets.Item("CreateProject")
$wsProject.Range("A3").Select()

# -----------------------------------------------------------------
# Sheet "createtask": customername values for rows 2-6 change from
# "cust2122" to "a212211" (row 7 is left as-is). This sheet becomes
# the active tab/sheet, so it is activated last.
# -----------------------------------------------------------------
$wsTask = $wb.Worksheets.Item("createtask")

$wsTask.Range("A2").Value = "a212211"
$wsTask.Range("A3").Value = "a212211"
$wsTask.Range("A4").Value = "a212211"
$wsTask.Range("A5").Value = "a212211"
$wsTask.Range("A6").Value = "a212211"

$wsTask.Activate()
$wsTask.Range("A2:A6").Select()
